$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 365 (this shifts rows 365:386 down to 366:387),
# copying formatting from the row above as Excel normally does.
$ws.Rows.Item(365).Insert()

# Populate the new row 365 with the new data record.
$ws.Cells.Item(365, 1).Value = 4
$ws.Cells.Item(365, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(365, 3).Value = "Los Lagos"
$ws.Cells.Item(365, 4).Value = 44610
$ws.Cells.Item(365, 5).Value = 10
$ws.Cells.Item(365, 6).Value = "Fruta"
$ws.Cells.Item(365, 7).Value = 100108
$ws.Cells.Item(365, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(365, 9).Value = 100108006
$ws.Cells.Item(365, 10).Value = "Plátano"
$ws.Cells.Item(365, 11).Value = "Sin especificar"
$ws.Cells.Item(365, 12).Value = "Primera Pintón"
$ws.Cells.Item(365, 13).Value = 1200
$ws.Cells.Item(365, 14).Value = 18000
$ws.Cells.Item(365, 15).Value = 20000
$ws.Cells.Item(365, 16).Value = 19000
$ws.Cells.Item(365, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(365, 18).Value = "Ecuador"
$ws.Cells.Item(365, 19).Value = 950
$ws.Cells.Item(365, 20).Value = 20

# Ensure the date cell keeps the expected date number format style applied to column D.
$ws.Cells.Item(365, 4).NumberFormat = $ws.Cells.Item(366, 4).NumberFormat()
